$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting the existing row 45 (and below) down to 46.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record.
$ws.Cells.Item(45, 1).Value = 8
$ws.Cells.Item(45, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(45, 3).Value = "Coquimbo"
$ws.Cells.Item(45, 4).Value = 44504
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(45, 6).Value = 100112052
$ws.Cells.Item(45, 7).Value = "Albahaca"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 760
$ws.Cells.Item(45, 11).Value = 3000
$ws.Cells.Item(45, 12).Value = 3500
$ws.Cells.Item(45, 13).Value = 3250
$ws.Cells.Item(45, 14).Value = "$/paquete"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 3250
$ws.Cells.Item(45, 17).Value = 1
$ws.Cells.Item(45, 18).Value = "Hortaliza"
